# Insert a new "setup" worksheet between "input_concentrations" and
# "component_name", holding a small Calorimeter / Initial volume table.

$wb = $excel.ActiveWorkbook

$before = $wb.Worksheets.Item("component_name")
$ws = $wb.Worksheets.Add($before)
$ws.Name = "setup"

$ws.Range("A1").Value = "Calorimeter"
$ws.Range("B1").Value = "DSC"
$ws.Range("A2").Value = "Initial volume"
$ws.Range("B2").Value = 15

$ws.Activate()
[void]$ws.Range("E15").Select()
